# Applies updated loading_percent results for the 380 kV case (Case_4_38).
# Only columns B, C, D, E, G, I, L, M change across data rows 2-25;
# columns F, H, J, K, N, O remain 0 and are left untouched, matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ 2 = 22.37372227330574; 3 = 16.16327233346146; 4 = 6.077913852366012; 5 = 7.526206967384444; 7 = 3.77653650844453; 9 = 51.7693511630784; 12 = 10.75681739180708; 13 = 18.9096567644111 }
    3 = @{ 2 = 22.28587230375038; 3 = 15.56972250320391; 4 = 5.976801115551976; 5 = 7.42870064576338; 7 = 3.782535546321491; 9 = 50.47576031848161; 12 = 10.75901058246742; 13 = 18.93712084349249 }
    4 = @{ 2 = 22.24376640731675; 3 = 15.19985834493625; 4 = 5.91602708738063; 5 = 7.367088767154117; 7 = 3.786395907730004; 9 = 49.66598219238842; 12 = 10.761953805982; 13 = 18.96004966829457 }
    5 = @{ 2 = 22.22958749001019; 3 = 15.04803411239212; 4 = 5.89162214729925; 5 = 7.34154900705817; 7 = 3.788013768470174; 9 = 49.33241383217405; 12 = 10.76355419205183; 13 = 18.97091172267248 }
    6 = @{ 2 = 22.22741304836941; 3 = 15.02276536592789; 4 = 5.887592438248122; 5 = 7.337282204912408; 7 = 3.788285121802635; 9 = 49.27681876573515; 12 = 10.76384413835275; 13 = 18.97280684937203 }
    7 = @{ 2 = 22.24356312048354; 3 = 15.19781488821149; 4 = 5.915696450955046; 5 = 7.366746072109843; 7 = 3.786417545325294; 9 = 49.66149761813019; 12 = 10.76197376641343; 13 = 18.96019001951029 }
    8 = @{ 2 = 22.3409793636916; 3 = 15.95989191155796; 4 = 6.042797560424603; 5 = 7.492951850141698; 7 = 3.778568401799459; 9 = 51.32672769567262; 12 = 10.75724196800594; 13 = 18.9178634660993 }
    9 = @{ 2 = 22.62547412140002; 3 = 17.40094607383447; 4 = 6.301018311958901; 5 = 7.726377077991578; 7 = 3.764569001277986; 9 = 54.45596685694408; 12 = 10.76065045551272; 13 = 18.88327106079666 }
    10 = @{ 2 = 22.89049649922995; 3 = 18.41490315491464; 4 = 6.494303391631452; 5 = 7.889008296111312; 7 = 3.75511708171024; 9 = 56.65564087220397; 12 = 10.77091265920829; 13 = 18.8877091942555 }
    11 = @{ 2 = 23.02290310230336; 3 = 18.86431353321538; 4 = 6.582612237254358; 5 = 7.961011923706958; 7 = 3.750994747449285; 9 = 57.63170675421896; 12 = 10.77726881153159; 13 = 18.89626804159801 }
    12 = @{ 2 = 23.07471124873447; 3 = 19.03264021762907; 4 = 6.616075399912225; 5 = 7.98798882530563; 7 = 3.74945896794457; 9 = 57.99755599720697; 12 = 10.77991843009742; 13 = 18.90045282241043 }
    13 = @{ 2 = 23.06347982695181; 3 = 18.99647278595284; 4 = 6.608868012848795; 5 = 7.982191808107668; 7 = 3.74978860612982; 9 = 57.91893477046342; 12 = 10.77933699556859; 13 = 18.89950953850473 }
    14 = @{ 2 = 23.02713216340462; 3 = 18.87819985603896; 4 = 6.58536500254389; 5 = 7.963237154548676; 7 = 3.750867893153852; 9 = 57.66188211675344; 12 = 10.77748193407719; 13 = 18.89659339632006 }
    15 = @{ 2 = 23.00508436404348; 3 = 18.80550847550415; 4 = 6.57097071766342; 5 = 7.951589066409483; 7 = 3.751532269646785; 9 = 57.50393300892618; 12 = 10.77637725704617; 13 = 18.89493016041953 }
    16 = @{ 2 = 22.88207892169974; 3 = 18.38528134969757; 4 = 6.488537096333278; 5 = 7.884262462580765; 7 = 3.755390034216771; 9 = 56.59133785946901; 12 = 10.77053123336417; 13 = 18.88728175415779 }
    17 = @{ 2 = 22.80962990091605; 3 = 18.12433762904989; 4 = 6.438041260168541; 5 = 7.842449573907162; 7 = 3.757801904937962; 9 = 56.02502431403568; 12 = 10.76737719605612; 13 = 18.88426734320211 }
    18 = @{ 2 = 22.76907659058643; 3 = 17.97314271132411; 4 = 6.409035275815368; 5 = 7.818213986892074; 7 = 3.759205861139395; 9 = 55.69699499317255; 12 = 10.76572195178021; 13 = 18.88314895685948 }
    19 = @{ 2 = 22.75553884741885; 3 = 17.92176554123771; 4 = 6.399221861922985; 5 = 7.809976454765738; 7 = 3.759684094452063; 9 = 55.58554213656485; 12 = 10.76518880034403; 13 = 18.88287587579545 }
    20 = @{ 2 = 22.81722680936352; 3 = 18.15223126701905; 4 = 6.443412945925544; 5 = 7.846919905371766; 7 = 3.757543429362696; 9 = 56.08554921860978; 12 = 10.76769650563973; 13 = 18.8845245119439 }
    21 = @{ 2 = 23.03776336702267; 3 = 18.91299095087086; 4 = 6.592268048569416; 5 = 7.968812484741833; 7 = 3.750550196868611; 9 = 57.7374886579701; 12 = 10.77802022503855; 13 = 18.89742430317574 }
    22 = @{ 2 = 23.19160362273157; 3 = 19.39932019561704; 4 = 6.689670486576996; 5 = 8.046787980537625; 7 = 3.746126826824578; 9 = 58.79508385345014; 12 = 10.78618181822202; 13 = 18.9113567936405 }
    23 = @{ 2 = 23.10862066210981; 3 = 19.1407971772909; 4 = 6.637684698164399; 5 = 8.005327000120161; 7 = 3.748474284902933; 9 = 58.23271359562787; 12 = 10.78169643695161; 13 = 18.90341647666991 }
    24 = @{ 2 = 22.81378882257745; 3 = 18.13962420108221; 4 = 6.4409843272443; 5 = 7.844899481133584; 7 = 3.757660232142116; 9 = 56.05819351680547; 12 = 10.76755165347399; 13 = 18.88440633175293 }
    25 = @{ 2 = 22.5385838527906; 3 = 17.01816923038606; 4 = 6.230398920367971; 5 = 7.664766191290556; 7 = 3.768208715951874; 9 = 53.62577442008024; 12 = 10.75836725090696; 13 = 18.88740840821069 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowData = $newValues[$rowNum]
    foreach ($colNum in $rowData.Keys) {
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$colNum]
    }
}

Write-Host "Updated loading_percent values for 380 kV case (rows 2-25)"